$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text, matching the source data which
# stores prices as literal strings (e.g. "30.705.40", "1.002") rather than numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.705.40"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "1.893.65"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "241.41"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "0.4912"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "0.2933"
$ws.Range("E8").Value = "  +0.74%  "
$ws.Range("D9").Value = "0.06742"
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("D10").Value = "1.893.35"
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("D11").Value = "17.15"
$ws.Range("E11").Value = "  +4.65%  "
$ws.Range("D12").Value = "0.07253"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "90.75"
$ws.Range("E13").Value = "  +5.09%  "
$ws.Range("D14").Value = "0.6750"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").Value = "5.029"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").Value = "30.647.07"
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("D17").Value = "0.000007985"
$ws.Range("E17").Value = "  +2.33%  "
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").Value = "13.12"
$ws.Range("E19").Value = "  +2.46%  "
$ws.Range("D20").Value = "2.140.04"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "4.799"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").Value = "191.07"
$ws.Range("E23").Value = "  +32.74%  "
$ws.Range("D24").Value = "6.089"
$ws.Range("E24").Value = "  +2.91%  "
$ws.Range("D25").Value = "9.376"
$ws.Range("E25").Value = "  +1.93%  "
$ws.Range("D26").Value = "156.88"
$ws.Range("E26").Value = "  +3.18%  "
$ws.Range("D27").Value = "18.88"
$ws.Range("E27").Value = "  +11.04%  "
$ws.Range("D28").Value = "1.891"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").Value = "1.408"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "4.293"
$ws.Range("E30").Value = "  +1.86%  "
$ws.Range("D31").Value = "0.09055"
$ws.Range("E31").Value = "  +2.74%  "
$ws.Range("D32").Value = "3.999"
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("D34").Value = "0.7403"
$ws.Range("E34").Value = "  +2.37%  "
$ws.Range("D35").Value = "1.104"
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("D36").Value = "2.753"
$ws.Range("E36").Value = "  +3.44%  "
$ws.Range("D37").Value = "0.01829"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("D39").Value = "0.9332"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").Value = "2.119"
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("D41").Value = "0.4404"
$ws.Range("E41").Value = "  +3.21%  "
$ws.Range("D42").Value = "105.02"
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "5.728"
$ws.Range("E44").Value = "  -0.66%  "
$ws.Range("D45").Value = "0.1353"
$ws.Range("E45").Value = "  +5.25%  "
$ws.Range("D46").Value = "7.516"
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("D47").Value = "0.05878"
$ws.Range("E47").Value = "  +2.50%  "
$ws.Range("D48").Value = "8.763"
$ws.Range("E48").Value = "  +5.30%  "
$ws.Range("D49").Value = "1.428"
$ws.Range("E49").Value = "  +5.44%  "
$ws.Range("D50").Value = "0.3956"
$ws.Range("E50").Value = "  +4.48%  "
$ws.Range("D51").Value = "33.82"
$ws.Range("E51").Value = "  +2.73%  "
